$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 odds update (Telstar - Venlo)
$ws.Range("N3").Value = 17      # Odd_Under05_FT: 15 -> 17
$ws.Range("O3").Value = 1.17    # Odd_Over15_FT: 1.18 -> 1.17
$ws.Range("P3").Value = 5       # Odd_Under15_FT: 4.5 -> 5
$ws.Range("Q3").Value = 1.57    # Odd_Over25_FT: 1.6 -> 1.57
$ws.Range("R3").Value = 2.35    # Odd_Under25_FT: 2.3 -> 2.35
